# v1.4.2 - Change web front end to use the Tombola RPM as the required function.
# Scale the Hz->RPM calibration table from 10..90 Hz to 1000..9000 "Hz" (tombola drum
# counts), flip the D column formula from B/A to A/B (RPM per unit), and add a
# calculator block (rows 14-18) that converts a target RPM (e.g. mains 60Hz =>
# 60 RPM) into the equivalent tombola rotation count using the averaged factor
# in D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the shared string used at C13/D14 ("RPM = " -> "1 RPM = ") ---
$ws.Range("C13").Value = "1 RPM = "

# --- Rescale the A column (Hz axis) from 10..90 to 1000..9000 ---
for ($i = 0; $i -lt 9; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = 1000 * ($i + 1)
}

# --- Flip the D column ratio formulas from B/A to A/B ---
$ws.Range("D3").Formula = "=A3/B3"
$ws.Range("D4:D11").Formula = "=A4/B4"

# --- New D12 cell: empty, but carrying the new 0.0000 number format ---
$ws.Range("D12").NumberFormat = "0.0000"

# --- New calculator block: rows 14-18 ---
$ws.Range("C14").Value = "RPM"
$ws.Range("D14").Value = "Hz"

$ws.Range("C15").Value = 59.9
$ws.Range("D15").Formula = "=INT(C15*`$D`$13)"
$ws.Range("E15").Formula = "=D15/`$D`$13"

$ws.Range("C16").Value = 60
$ws.Range("D16").Formula = "=INT(C16*`$D`$13)"

$ws.Range("C17").Value = 60.1
$ws.Range("D17").Formula = "=INT(C17*`$D`$13)"

$ws.Range("C18").Value = 0.1
$ws.Range("D18").Formula = "=INT(C18*`$D`$13)"

# E16:E18 is one shared formula (D/$D$13) filled down from E16
$ws.Range("E16:E18").Formula = "=D16/`$D`$13"

# --- Column D width (autofit-ish, matches the bestFit width in the target) ---
$ws.Columns.Item(4).ColumnWidth = 9.36328125

# --- Selection matches the post-edit cursor position ---
$ws.Range("D18").Select()

$wb.Application.Calculate()
